$d = $word.ActiveDocument

# Locate the end of the last existing HULA entry (HULA-026's "Para:" line)
# so that the new HULA-027 user story can be appended right after it, in the
# same place a human editor would have placed it (before the existing blank
# separator paragraphs that lead into the document's closing section).
$anchorText = "Para: Mantener la información actualizada para una mejor gestión del taller y sus recursos"

$searchRange = $d.Content
$found = $searchRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph for HULA-026 'Para:' line."
}

$insertionPoint = $searchRange.Duplicate
$insertionPoint.Collapse(0)

$cr = [string][char]13

$newContent = $cr + $cr + `
    "HULA-027 Consulta de estado de carro en el Taller" + $cr + `
    "Como: Empleado del área" + $cr + `
    "Quiero: Consultar el estado del carro en el taller (estadocarro) y los repuestos disponibles" + $cr + `
    "Para: Informar al cliente sobre el avance de la reparación y verificar disponibilidad de repuestos"

$insertionPoint.InsertAfter($newContent)
